# cryptos.xlsx refresh -- GitHub Actions scheduled data pull.
# Rewrites the Price (D) / Volume(1h) (E) figures for each coin row with
# the latest scrape, and (for rows 31-32) the coin Name/Link/Price too,
# since NEARProtocol and EthereumClassic swapped rank that run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
    $ws.Range("D2").Value = "59.885.97"
    $ws.Range("E2").Value = "  +3.75%  "
# Row 3
    $ws.Range("D3").Value = "3.020.86"
    $ws.Range("E3").Value = "  +2.82%  "
# Row 4
    $ws.Range("E4").Value = "  +0.17%  "
# Row 5
    $ws.Range("D5").Value = "'565.87"
    $ws.Range("D5").Style = "Normal"
    $ws.Range("E5").Value = "  +2.90%  "
# Row 6
    $ws.Range("D6").Value = "'140.72"
    $ws.Range("D6").Style = "Normal"
    $ws.Range("E6").Value = "  +7.56%  "
# Row 7
    $ws.Range("E7").Value = "  -0.04%  "
# Row 8
    $ws.Range("E8").Value = "  +1.97%  "
# Row 9
    $ws.Range("D9").Value = "3.011.88"
    $ws.Range("E9").Value = "  +2.86%  "
# Row 11
    $ws.Range("D11").Value = "'5.32"
    $ws.Range("D11").Style = "Normal"
    $ws.Range("E11").Value = "  +11.25%  "
# Row 12
    $ws.Range("D12").Value = "'0.462"
    $ws.Range("D12").Style = "Normal"
    $ws.Range("E12").Value = "  +3.49%  "
# Row 13
    $ws.Range("E13").Value = "  +5.24%  "
# Row 14
    $ws.Range("D14").Value = "'34.28"
    $ws.Range("D14").Style = "Normal"
    $ws.Range("E14").Value = "  +4.04%  "
# Row 15
    $ws.Range("E15").Value = "  +1.85%  "
# Row 16
    $ws.Range("D16").Value = "3.520.19"
    $ws.Range("E16").Value = "  +2.88%  "
# Row 17
    $ws.Range("D17").Value = "'7.22"
    $ws.Range("D17").Style = "Normal"
    $ws.Range("E17").Value = "  +5.37%  "
# Row 18
    $ws.Range("D18").Value = "3.021.74"
    $ws.Range("E18").Value = "  +3.15%  "
# Row 19
    $ws.Range("D19").Value = "59.929.64"
    $ws.Range("E19").Value = "  +3.90%  "
# Row 20
    $ws.Range("D20").Value = "'438.60"
    $ws.Range("D20").Style = "Normal"
    $ws.Range("E20").Value = "  +4.84%  "
# Row 21
    $ws.Range("D21").Value = "'13.73"
    $ws.Range("D21").Style = "Normal"
    $ws.Range("E21").Value = "  +4.15%  "
# Row 22
    $ws.Range("D22").Value = "'0.724"
    $ws.Range("D22").Style = "Normal"
    $ws.Range("E22").Value = "  +5.48%  "
# Row 23
    $ws.Range("D23").Value = "'7.15"
    $ws.Range("D23").Style = "Normal"
    $ws.Range("E23").Value = "  +2.39%  "
# Row 24
    $ws.Range("D24").Value = "'13.31"
    $ws.Range("D24").Style = "Normal"
    $ws.Range("E24").Value = "  +2.00%  "
# Row 25
    $ws.Range("D25").Value = "'80.87"
    $ws.Range("D25").Style = "Normal"
    $ws.Range("E25").Value = "  +1.30%  "
# Row 26
    $ws.Range("E26").Value = "  -0.10%  "
# Row 27
    $ws.Range("D27").Value = "'2.27"
    $ws.Range("D27").Style = "Normal"
    $ws.Range("E27").Value = "  +14.17%  "
# Row 28
    $ws.Range("E28").Value = "  +0.51%  "
# Row 29
    $ws.Range("E29").Value = "  +3.48%  "
# Row 30
    $ws.Range("D30").Value = "'7.88"
    $ws.Range("D30").Style = "Normal"
    $ws.Range("E30").Value = "  +5.14%  "
# Row 31
    $ws.Range("B31").Value = "EthereumClassic"
    $ws.Range("C31").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
    $ws.Range("D31").Value = "'26.13"
    $ws.Range("D31").Style = "Normal"
    $ws.Range("E31").Value = "  +3.82%  "
# Row 32
    $ws.Range("B32").Value = "NEARProtocol"
    $ws.Range("C32").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
    $ws.Range("D32").Value = "'6.31"
    $ws.Range("D32").Style = "Normal"
    $ws.Range("E32").Value = "  +5.14%  "
# Row 33
    $ws.Range("E33").Value = "  +4.76%  "
# Row 34
    $ws.Range("E34").Value = "  +15.80%  "
# Row 35
    $ws.Range("D35").Value = "'1.00"
    $ws.Range("D35").Style = "Normal"
    $ws.Range("E35").Value = "  +6.92%  "
# Row 36
    $ws.Range("E36").Value = "  +5.03%  "
# Row 37
    $ws.Range("D37").Value = "'2.12"
    $ws.Range("D37").Style = "Normal"
    $ws.Range("E37").Value = "  +2.14%  "
# Row 38
    $ws.Range("E38").Value = "  +2.56%  "
# Row 39
    $ws.Range("D39").Value = "'8.68"
    $ws.Range("D39").Style = "Normal"
    $ws.Range("E39").Value = "  -0.49%  "
# Row 40
    $ws.Range("E40").Value = "  +9.09%  "
# Row 41
    $ws.Range("D41").Value = "'406.70"
    $ws.Range("D41").Style = "Normal"
    $ws.Range("E41").Value = "  +7.97%  "
# Row 42
    $ws.Range("D42").Value = "'0.0355"
    $ws.Range("D42").Style = "Normal"
    $ws.Range("E42").Value = "  +2.64%  "
# Row 43
    $ws.Range("D43").Value = "2.786.37"
    $ws.Range("E43").Value = "  +3.87%  "
# Row 44
    $ws.Range("E44").Value = "  -0.33%  "
# Row 45
    $ws.Range("E45").Value = "  +6.48%  "
# Row 46
    $ws.Range("E46").Value = "  +0.01%  "
# Row 47
    $ws.Range("D47").Value = "'123.24"
    $ws.Range("D47").Style = "Normal"
    $ws.Range("E47").Value = "  +0.93%  "
# Row 48
    $ws.Range("E48").Value = "  +3.21%  "
# Row 49
    $ws.Range("E49").Value = "  +1.65%  "
# Row 50
    $ws.Range("D50").Value = "'34.13"
    $ws.Range("D50").Style = "Normal"
    $ws.Range("E50").Value = "  +20.24%  "
# Row 51
    $ws.Range("D51").Value = "'23.68"
    $ws.Range("D51").Style = "Normal"
    $ws.Range("E51").Value = "  +2.10%  "
